# Generate Report for Handback
#
# Fills in the "Latest Target File", "Latest Handback File",
# "Latest Handback DateTime" and "Error Detail" columns (I, J, K, P) for
# row 7 (the 0945156b-8fed-4cb8-95c3-9698feda23e8 entry) on both the
# zh-cn and de-de worksheets, because a handback was received for that
# file but its version did not match the latest handoff.

$wb = $excel.ActiveWorkbook

$latestFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6904de89948fe05c2cd06a4968960cab1113d343/e2e/0945156b-8fed-4cb8-95c3-9698feda23e8.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99f1a4ade8cba079817d19868b19611ad6f8c93f/e2e/0945156b-8fed-4cb8-95c3-9698feda23e8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6904de89948fe05c2cd06a4968960cab1113d343/e2e/0945156b-8fed-4cb8-95c3-9698feda23e8.md."

# --- zh-cn worksheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestFileUrl, [System.Type]::Missing, [System.Type]::Missing, "0945156b-8fed-4cb8-95c3-9698feda23e8.md")
$wsZh.Range("J7").Value = "0945156b-8fed-4cb8-95c3-9698feda23e8.937881b94b40b501fcdabce08ec119eb7886dfe2.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-26 02:55:18"
$wsZh.Range("P7").Value = $errorDetail

# --- de-de worksheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestFileUrl, [System.Type]::Missing, [System.Type]::Missing, "0945156b-8fed-4cb8-95c3-9698feda23e8.md")
$wsDe.Range("J7").Value = "0945156b-8fed-4cb8-95c3-9698feda23e8.937881b94b40b501fcdabce08ec119eb7886dfe2.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-26 02:55:24"
$wsDe.Range("P7").Value = $errorDetail
